$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Flt3l"
$ws.Range("C2").Value = "Flt3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 7.556524
$ws.Range("H2").Value = 22.669572
$ws.Range("I2").Value = 0.3236988808488902
$ws.Range("J2").Value = 0.3236988808488901
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 1.353118333333333
$ws.Range("N2").Value = 4.059355
$ws.Range("O2").Value = 0.5236944229679847
$ws.Range("P2").Value = 0.5236944229679846
$ws.Range("Q2").Value = 10.22487116067333
$ws.Range("R2").Value = 92.02384044606
$ws.Range("S2").Value = 0.169519298621542
$ws.Range("T2").Value = 0.1695192986215419

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Flt3l"
$ws.Range("C3").Value = "Flt3"
$ws.Range("D3").Value = "M2"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 7.556524
$ws.Range("H3").Value = 22.669572
$ws.Range("I3").Value = 0.3236988808488902
$ws.Range("J3").Value = 0.3236988808488901
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.230675333333333
$ws.Range("N3").Value = 3.692026
$ws.Range("O3").Value = 0.4763055770320153
$ws.Range("P3").Value = 0.4763055770320153
$ws.Range("Q3").Value = 9.299627692541334
$ws.Range("R3").Value = 83.69664923287202
$ws.Range("S3").Value = 0.1541795822273482
$ws.Range("T3").Value = 0.1541795822273482

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Flt3l"
$ws.Range("C4").Value = "Flt3"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 7.625735333333334
$ws.Range("H4").Value = 22.877206
$ws.Range("I4").Value = 0.326663687305147
$ws.Range("J4").Value = 0.3266636873051469
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 1.353118333333333
$ws.Range("N4").Value = 4.059355
$ws.Range("O4").Value = 0.5236944229679847
$ws.Range("P4").Value = 0.5236944229679846
$ws.Range("Q4").Value = 10.31852228468111
$ws.Range("R4").Value = 92.86670056213
$ws.Range("S4").Value = 0.1710719512278632
$ws.Range("T4").Value = 0.1710719512278631

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Flt3l"
$ws.Range("C5").Value = "Flt3"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 7.625735333333334
$ws.Range("H5").Value = 22.877206
$ws.Range("I5").Value = 0.326663687305147
$ws.Range("J5").Value = 0.3266636873051469
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.230675333333333
$ws.Range("N5").Value = 3.692026
$ws.Range("O5").Value = 0.4763055770320153
$ws.Range("P5").Value = 0.4763055770320153
$ws.Range("Q5").Value = 9.384804373261778
$ws.Range("R5").Value = 84.46323935935601
$ws.Range("S5").Value = 0.1555917360772839
$ws.Range("T5").Value = 0.1555917360772838

# Row 6
$ws.Range("A6").Value = "M2"
$ws.Range("B6").Value = "Flt3l"
$ws.Range("C6").Value = "Flt3"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 5.476188666666666
$ws.Range("H6").Value = 16.428566
$ws.Range("I6").Value = 0.2345835390342671
$ws.Range("J6").Value = 0.234583539034267
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 1.353118333333333
$ws.Range("N6").Value = 4.059355
$ws.Range("O6").Value = 0.5236944229679847
$ws.Range("P6").Value = 0.5236944229679846
$ws.Range("Q6").Value = 7.409931281658889
$ws.Range("R6").Value = 66.68938153493
$ws.Range("S6").Value = 0.1228500911123382
$ws.Range("T6").Value = 0.1228500911123382

# Row 7
$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "Flt3l"
$ws.Range("C7").Value = "Flt3"
$ws.Range("D7").Value = "M2"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 5.476188666666666
$ws.Range("H7").Value = 16.428566
$ws.Range("I7").Value = 0.2345835390342671
$ws.Range("J7").Value = 0.234583539034267
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.230675333333333
$ws.Range("N7").Value = 3.692026
$ws.Range("O7").Value = 0.4763055770320153
$ws.Range("P7").Value = 0.4763055770320153
$ws.Range("Q7").Value = 6.739410312746222
$ws.Range("R7").Value = 60.65469281471601
$ws.Range("S7").Value = 0.1117334479219289
$ws.Range("T7").Value = 0.1117334479219289

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Flt3l"
$ws.Range("C8").Value = "Flt3"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2.685852666666667
$ws.Range("H8").Value = 8.057558
$ws.Range("I8").Value = 0.1150538928116959
$ws.Range("J8").Value = 0.1150538928116958
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 1.353118333333333
$ws.Range("N8").Value = 4.059355
$ws.Range("O8").Value = 0.5236944229679847
$ws.Range("P8").Value = 0.5236944229679846
$ws.Range("Q8").Value = 3.634276483898889
$ws.Range("R8").Value = 32.70848835509
$ws.Range("S8").Value = 0.06025308200624144
$ws.Range("T8").Value = 0.06025308200624141

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Flt3l"
$ws.Range("C9").Value = "Flt3"
$ws.Range("D9").Value = "M2"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2.685852666666667
$ws.Range("H9").Value = 8.057558
$ws.Range("I9").Value = 0.1150538928116959
$ws.Range("J9").Value = 0.1150538928116958
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.230675333333333
$ws.Range("N9").Value = 3.692026
$ws.Range("O9").Value = 0.4763055770320153
$ws.Range("P9").Value = 0.4763055770320153
$ws.Range("Q9").Value = 3.305412625834223
$ws.Range("R9").Value = 29.748713632508
$ws.Range("S9").Value = 0.05480081080545444
$ws.Range("T9").Value = 0.05480081080545442
